# Configsheet.xlsx - "Commiting with Configuration Definition"
# Add a configuration-definition index column (C) to the "Bill Configuration"
# sheet: rows 2-17 get sequential values 1-16, mirroring the existing
# configuration rows in columns A/B. Also updates the active selection on
# that sheet to J13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill Configuration")
$ws.Activate()

for ($i = 0; $i -lt 16; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 3).Value = $i + 1
}

$ws.Range("J13").Select()
